# Auto-generated Excel COM-interop script
# Applies numeric value updates to Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled Universalis market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2918.125
$ws.Range("I17").Value = 2957.5
$ws.Range("J17").Value = 2800
$ws.Range("K17").Value = 8872.5
$ws.Range("L17").Value = 8400
$ws.Range("M17").Value = -8704.5
$ws.Range("N17").Value = -8736

# Row 39
$ws.Range("H39").Value = 1096.6
$ws.Range("I39").Value = 1096.6
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3289.8
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2993.8
$ws.Range("N39").Value = $null

# Row 74
$ws.Range("H74").Value = 1659.8
$ws.Range("I74").Value = 1659.8
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1659.8
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -723.8
$ws.Range("N74").Value = $null

# Row 77
$ws.Range("H77").Value = 1659.8
$ws.Range("I77").Value = 1659.8
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8299
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3619
$ws.Range("N77").Value = $null

# Row 116
$ws.Range("H116").Value = 9174.75
$ws.Range("I116").Value = 9174.75
$ws.Range("K116").Value = 9174.75
$ws.Range("M116").Value = -5732.75

# Row 138
$ws.Range("H138").Value = 2463.9412
$ws.Range("J138").Value = 3193.8
$ws.Range("L138").Value = 9581.400000000001
$ws.Range("N138").Value = -19861.4


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1327.75
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 1768.6666
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1768.6666
$ws.Range("M3").Value = 110
$ws.Range("N3").Value = -1998.6666

# Row 17
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -2346

# Row 22
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 10000
$ws.Range("K22").Value = 10000
$ws.Range("M22").Value = -9701

# Row 97
$ws.Range("H97").Value = 166669170
$ws.Range("I97").Value = 166669170
$ws.Range("K97").Value = 166669170
$ws.Range("M97").Value = -166668674

# Row 102
$ws.Range("H102").Value = 17585318
$ws.Range("I102").Value = 1376752.4
$ws.Range("K102").Value = 1376752.4
$ws.Range("M102").Value = -1375130.4


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 1067
$ws.Range("I25").Value = 750.5
$ws.Range("J25").Value = 1700
$ws.Range("K25").Value = 750.5
$ws.Range("L25").Value = 1700
$ws.Range("M25").Value = -515.5
$ws.Range("N25").Value = -2170

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null

# Row 86
$ws.Range("H86").Value = 4133
$ws.Range("I86").Value = 4133
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4133
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3010
$ws.Range("N86").Value = $null

# Row 89
$ws.Range("H89").Value = 4133
$ws.Range("I89").Value = 4133
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 20665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -15049
$ws.Range("N89").Value = $null

# Row 99
$ws.Range("H99").Value = 2345.0588
$ws.Range("I99").Value = 2345.0588
$ws.Range("K99").Value = 2345.0588
$ws.Range("M99").Value = -847.0587999999998

# Row 105
$ws.Range("H105").Value = 2002.3334
$ws.Range("I105").Value = 2001.5
$ws.Range("J105").Value = 2009
$ws.Range("K105").Value = 2001.5
$ws.Range("L105").Value = 2009
$ws.Range("M105").Value = -254.5
$ws.Range("N105").Value = -5503

# Row 107
$ws.Range("H107").Value = 47816.11
$ws.Range("I107").Value = 53293.125
$ws.Range("K107").Value = 53293.125
$ws.Range("M107").Value = -51373.125


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 612.2
$ws.Range("J22").Value = 416.85715
$ws.Range("L22").Value = 416.85715
$ws.Range("N22").Value = -1116.85715

# Row 107
$ws.Range("H107").Value = 1033.3334
$ws.Range("J107").Value = 1050
$ws.Range("L107").Value = 1050
$ws.Range("N107").Value = -4890


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 350.16666
$ws.Range("I2").Value = 499.66666
$ws.Range("J2").Value = 200.66667
$ws.Range("K2").Value = 2997.99996
$ws.Range("L2").Value = 1204.00002
$ws.Range("M2").Value = -2884.99996
$ws.Range("N2").Value = -1430.00002

# Row 12
$ws.Range("H12").Value = 210.375
$ws.Range("J12").Value = 305.54544
$ws.Range("L12").Value = 916.63632
$ws.Range("N12").Value = -1262.63632

# Row 100
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null

# Row 134
$ws.Range("H134").Value = 3378
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

# Row 141
$ws.Range("H141").Value = 1373.5
$ws.Range("I141").Value = 1182.1666
$ws.Range("J141").Value = 1947.5
$ws.Range("K141").Value = 3546.4998
$ws.Range("L141").Value = 5842.5
$ws.Range("M141").Value = 1633.5002
$ws.Range("N141").Value = -16202.5


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5164.6665
$ws.Range("I80").Value = 5247
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 5247
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -4249
$ws.Range("N80").Value = -6996

# Row 83
$ws.Range("H83").Value = 5164.6665
$ws.Range("I83").Value = 5247
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 26235
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -21243
$ws.Range("N83").Value = -34984

# Row 97
$ws.Range("H97").Value = 469.86667
$ws.Range("I97").Value = 218
$ws.Range("K97").Value = 218
$ws.Range("M97").Value = 278

# Row 102
$ws.Range("H102").Value = 1915.5714
$ws.Range("I102").Value = 1984.8334
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1984.8334
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -362.8334
$ws.Range("N102").Value = -4744

# Row 113
$ws.Range("H113").Value = 525
$ws.Range("J113").Value = 300
$ws.Range("L113").Value = 300
$ws.Range("N113").Value = -4640


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1074.1111
$ws.Range("J22").Value = 1277.7778
$ws.Range("L22").Value = 1277.7778
$ws.Range("N22").Value = -1867.7778

# Row 27
$ws.Range("H27").Value = 1074.1111
$ws.Range("J27").Value = 1277.7778
$ws.Range("L27").Value = 1277.7778
$ws.Range("N27").Value = -1491.7778

# Row 46
$ws.Range("H46").Value = 933
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 966
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 966
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1342

# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = $null

# Row 61
$ws.Range("H61").Value = 3943
$ws.Range("I61").Value = 3433.5
$ws.Range("K61").Value = 3433.5
$ws.Range("M61").Value = -3231.5

# Row 68
$ws.Range("H68").Value = 3505.182
$ws.Range("J68").Value = 3796.6
$ws.Range("L68").Value = 3796.6
$ws.Range("N68").Value = -5294.6

# Row 71
$ws.Range("H71").Value = 3505.182
$ws.Range("J71").Value = 3796.6
$ws.Range("L71").Value = 18983
$ws.Range("N71").Value = -26471

# Row 100
$ws.Range("H100").Value = 2307.25
$ws.Range("I100").Value = 1818.7
$ws.Range("J100").Value = 4750
$ws.Range("K100").Value = 1818.7
$ws.Range("L100").Value = 4750
$ws.Range("M100").Value = -1277.7
$ws.Range("N100").Value = -5832

# Row 113
$ws.Range("H113").Value = 3943
$ws.Range("I113").Value = 3433.5
$ws.Range("K113").Value = 3433.5
$ws.Range("M113").Value = -1263.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 62626
$ws.Range("J45").Value = 62626
$ws.Range("L45").Value = 62626
$ws.Range("N45").Value = -63608

# Row 62
$ws.Range("H62").Value = 6001
$ws.Range("I62").Value = 6001
$ws.Range("K62").Value = 6001
$ws.Range("M62").Value = -5377

# Row 65
$ws.Range("H65").Value = 6001
$ws.Range("I65").Value = 6001
$ws.Range("K65").Value = 30005
$ws.Range("M65").Value = -26885

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = $null

# Row 126
$ws.Range("H126").Value = 4708.727
$ws.Range("I126").Value = 3775.25
$ws.Range("K126").Value = 11325.75
$ws.Range("M126").Value = -8855.75

